$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 4650.3438
$ws.Range("I40").Value = 2595.3333
$ws.Range("K40").Value = 2595.3333
$ws.Range("M40").Value = -2420.3333
$ws.Range("H113").Value = 62523012
$ws.Range("I113").Value = 83337930
$ws.Range("K113").Value = 83337930
$ws.Range("M113").Value = -83334676
$ws.Range("H135").Value = 829.9
$ws.Range("I135").Value = 829.9
$ws.Range("K135").Value = 7469.099999999999
$ws.Range("M135").Value = -4934.099999999999
$ws.Range("H138").Value = 3023.5715
$ws.Range("J138").Value = 3253.5715
$ws.Range("L138").Value = 9760.7145
$ws.Range("N138").Value = -20040.7145

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1609.4615
$ws.Range("I2").Value = 1190.0454
$ws.Range("K2").Value = 1190.0454
$ws.Range("M2").Value = -1077.0454
$ws.Range("H32").Value = 8788.757
$ws.Range("I32").Value = 8866.25
$ws.Range("K32").Value = 8866.25
$ws.Range("M32").Value = -8579.25
$ws.Range("H45").Value = 5840.5625
$ws.Range("J45").Value = 7900
$ws.Range("L45").Value = 7900
$ws.Range("N45").Value = -8654
$ws.Range("H74").Value = 4289.5356
$ws.Range("I74").Value = 3324.48
$ws.Range("K74").Value = 3324.48
$ws.Range("M74").Value = -2450.48
$ws.Range("H77").Value = 4289.5356
$ws.Range("I77").Value = 3324.48
$ws.Range("K77").Value = 16622.4
$ws.Range("M77").Value = -12254.4
$ws.Range("H116").Value = 1609.4615
$ws.Range("I116").Value = 1190.0454
$ws.Range("K116").Value = 1190.0454
$ws.Range("M116").Value = 1103.9546
$ws.Range("H132").Value = 1389
$ws.Range("I132").Value = 1286.3334
$ws.Range("K132").Value = 3859.0002
$ws.Range("M132").Value = -1329.0002
$ws.Range("H135").Value = 59648.93
$ws.Range("J135").Value = 59648.93
$ws.Range("L135").Value = 59648.93
$ws.Range("N135").Value = -69788.92999999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1609.4615
$ws.Range("I3").Value = 1190.0454
$ws.Range("K3").Value = 1190.0454
$ws.Range("M3").Value = -1076.0454
$ws.Range("H86").Value = 450
$ws.Range("I86").Value = 450
$ws.Range("K86").Value = 450
$ws.Range("M86").Value = 673
$ws.Range("H89").Value = 450
$ws.Range("I89").Value = 450
$ws.Range("K89").Value = 2250
$ws.Range("M89").Value = 3366
$ws.Range("H99").Value = 3828.8076
$ws.Range("I99").Value = 2934.625
$ws.Range("K99").Value = 2934.625
$ws.Range("M99").Value = -1436.625
$ws.Range("H105").Value = 2730.6
$ws.Range("J105").Value = 2587.4
$ws.Range("L105").Value = 2587.4
$ws.Range("N105").Value = -6081.4
$ws.Range("H106").Value = 198333.33
$ws.Range("I106").Value = 500000
$ws.Range("J106").Value = 47500
$ws.Range("K106").Value = 500000
$ws.Range("L106").Value = 47500
$ws.Range("M106").Value = -498738
$ws.Range("N106").Value = -50024
$ws.Range("H134").Value = 4441.3125
$ws.Range("I134").Value = 4107.6553
$ws.Range("K134").Value = 12322.9659
$ws.Range("M134").Value = -9787.965900000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 1031.8462
$ws.Range("I22").Value = 1031.8462
$ws.Range("K22").Value = 1031.8462
$ws.Range("M22").Value = -681.8462
$ws.Range("H31").Value = 41673556
$ws.Range("I31").Value = 76926770
$ws.Range("J31").Value = 10672.637
$ws.Range("K31").Value = 76926770
$ws.Range("L31").Value = 10672.637
$ws.Range("M31").Value = -76926475
$ws.Range("N31").Value = -11262.637
$ws.Range("H34").Value = 41673556
$ws.Range("I34").Value = 76926770
$ws.Range("J34").Value = 10672.637
$ws.Range("K34").Value = 76926770
$ws.Range("L34").Value = 10672.637
$ws.Range("M34").Value = -76926568
$ws.Range("N34").Value = -11076.637
$ws.Range("H58").Value = 10349.333
$ws.Range("I58").Value = 7263.3335
$ws.Range("K58").Value = 7263.3335
$ws.Range("M58").Value = -7060.3335
$ws.Range("H62").Value = 5375
$ws.Range("H65").Value = 5375
$ws.Range("H136").Value = 10349.333
$ws.Range("I136").Value = 7263.3335
$ws.Range("K136").Value = 21790.0005
$ws.Range("M136").Value = -19240.0005

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H50").Value = 877.7857
$ws.Range("I50").Value = 1743
$ws.Range("K50").Value = 5229
$ws.Range("M50").Value = -4748
$ws.Range("H53").Value = 877.7857
$ws.Range("I53").Value = 1743
$ws.Range("K53").Value = 5229
$ws.Range("M53").Value = -4748
$ws.Range("H107").Value = 312800.84
$ws.Range("I107").Value = 264.66666
$ws.Range("K107").Value = 793.9999799999999
$ws.Range("M107").Value = 1126.00002
$ws.Range("H113").Value = 1988.8334
$ws.Range("J113").Value = 2249.7
$ws.Range("L113").Value = 6749.099999999999
$ws.Range("N113").Value = -11089.1
$ws.Range("H131").Value = 23812612
$ws.Range("J131").Value = 3368.8333
$ws.Range("L131").Value = 10106.4999
$ws.Range("N131").Value = -20186.4999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H43").Value = 583
$ws.Range("I43").Value = 583
$ws.Range("K43").Value = 583
$ws.Range("M43").Value = -432

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 13725.117
$ws.Range("I93").Value = 2128.4167
$ws.Range("J93").Value = 41557.2
$ws.Range("K93").Value = 2128.4167
$ws.Range("L93").Value = 41557.2
$ws.Range("M93").Value = -880.4167000000002
$ws.Range("N93").Value = -44053.2
$ws.Range("H100").Value = 4654.0557
$ws.Range("I100").Value = 3779.3635
$ws.Range("J100").Value = 6028.5713
$ws.Range("K100").Value = 3779.3635
$ws.Range("L100").Value = 6028.5713
$ws.Range("M100").Value = -3238.3635
$ws.Range("N100").Value = -7110.5713

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H2").Value = 1722.9231
$ws.Range("I2").Value = 924.75
$ws.Range("K2").Value = 924.75
$ws.Range("M2").Value = -812.75
$ws.Range("H4").Value = 11808.714
$ws.Range("I4").Value = 11808.714
$ws.Range("K4").Value = 11808.714
$ws.Range("M4").Value = -11695.714
$ws.Range("H15").Value = 7267.5
$ws.Range("I15").Value = 0
$ws.Range("K15").Value = 0
$ws.Range("H45").Value = 27320.715
$ws.Range("J45").Value = 30249.8
$ws.Range("L45").Value = 30249.8
$ws.Range("N45").Value = -31231.8
$ws.Range("H100").Value = 866.5789
$ws.Range("I100").Value = 776.63635
$ws.Range("J100").Value = 990.25
$ws.Range("K100").Value = 1553.2727
$ws.Range("L100").Value = 1980.5
$ws.Range("M100").Value = -1012.2727
$ws.Range("N100").Value = -3062.5
$ws.Range("H136").Value = 2305.4666
$ws.Range("I136").Value = 1469.037
$ws.Range("K136").Value = 4407.111
$ws.Range("M136").Value = -1857.111

# Remove M15 cell entirely in WVR sheet (column M no longer present for this row)
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("M15").ClearContents()
